# Auto-generated edit script: refresh crypto price/volume/hour snapshot values
# Mirrors the "Updated symbol list on Mon Jan  2 10:13:27 UTC 2023 with GitHub
# Actions" commit: price (D), 1h volume % (E) and hour (G) columns are
# refreshed for every coin row; the date column (F) is unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $text) {
    # The sheet stores every data cell as literal text (inline string),
    # including values that look numeric/percentage ("246.55", "0.95%", "9").
    # Assigning such a string straight to .Value lets Excel auto-convert it
    # to a real number, so prefix with an apostrophe (forces text entry),
    # then strip the resulting "Text" number-format/quote-prefix style back
    # to Normal so only the cell's value - not its formatting - changes.
    $ws.Range($ref).Value = "'" + $text
    $ws.Range($ref).Style = "Normal"
}

# Row 2
Set-TextCell "D2" "246.91"
Set-TextCell "E2" "1.17%"
Set-TextCell "G2" "10"

# Row 3
Set-TextCell "D3" "30.19"
Set-TextCell "E3" "11.15%"
Set-TextCell "G3" "10"

# Row 4
Set-TextCell "D4" "5.182"
Set-TextCell "E4" "0.69%"
Set-TextCell "G4" "10"

# Row 5
Set-TextCell "D5" "0.05740"
Set-TextCell "E5" "2.11%"
Set-TextCell "G5" "10"

# Row 6
Set-TextCell "D6" "6.594"
Set-TextCell "E6" "1.86%"
Set-TextCell "G6" "10"

# Row 7
Set-TextCell "D7" "0.8588"
Set-TextCell "E7" "5.18%"
Set-TextCell "G7" "10"

# Row 8
Set-TextCell "D8" "0.8748"
Set-TextCell "E8" "5.22%"
Set-TextCell "G8" "10"

# Row 9
Set-TextCell "D9" "0.1360"
Set-TextCell "E9" "2.47%"
Set-TextCell "G9" "10"

# Row 10
Set-TextCell "E10" "2.49%"
Set-TextCell "G10" "10"

# Row 11
Set-TextCell "D11" "0.02919"
Set-TextCell "E11" "1.22%"
Set-TextCell "G11" "10"

# Row 12
Set-TextCell "D12" "0.09395"
Set-TextCell "E12" "0.14%"
Set-TextCell "G12" "10"

# Row 13
Set-TextCell "D13" "0.001513"
Set-TextCell "E13" "0.27%"
Set-TextCell "G13" "10"

# Row 14
Set-TextCell "D14" "0.04146"
Set-TextCell "E14" "-1.29%"
Set-TextCell "G14" "10"

# Row 15
Set-TextCell "D15" "0.0006015"
Set-TextCell "E15" "-93.99%"
Set-TextCell "G15" "10"

# Row 16
Set-TextCell "D16" "0.006073"
Set-TextCell "E16" "-1.67%"
Set-TextCell "G16" "10"

# Row 17
Set-TextCell "D17" "3.505"
Set-TextCell "E17" "-2.82%"
Set-TextCell "G17" "10"

# Row 18
Set-TextCell "D18" "3.052"
Set-TextCell "E18" "1.01%"
Set-TextCell "G18" "10"

# Row 19
Set-TextCell "D19" "2.281"
Set-TextCell "E19" "-1.11%"
Set-TextCell "G19" "10"

# Row 20
Set-TextCell "G20" "10"

# Row 21
Set-TextCell "D21" "0.03280"
Set-TextCell "E21" "6.11%"
Set-TextCell "G21" "10"

# Row 22
Set-TextCell "D22" "0.1308"
Set-TextCell "E22" "1.30%"
Set-TextCell "G22" "10"

# Row 23
Set-TextCell "D23" "3.606"
Set-TextCell "E23" "-3.63%"
Set-TextCell "G23" "10"

# Row 24
Set-TextCell "E24" "0.46%"
Set-TextCell "G24" "10"

# Row 25
Set-TextCell "D25" "0.001215"
Set-TextCell "E25" "-0.77%"
Set-TextCell "G25" "10"

# Row 26
Set-TextCell "D26" "0.004505"
Set-TextCell "E26" "0.37%"
Set-TextCell "G26" "10"

# Row 27
Set-TextCell "D27" "0.0001179"
Set-TextCell "E27" "20.35%"
Set-TextCell "G27" "10"

# Row 28
Set-TextCell "E28" "0.54%"
Set-TextCell "G28" "10"

# Row 29
Set-TextCell "G29" "10"

# Row 30
Set-TextCell "G30" "10"

# Row 31
Set-TextCell "G31" "10"

# Row 32
Set-TextCell "G32" "10"

# Row 33
Set-TextCell "G33" "10"

# Row 34
Set-TextCell "G34" "10"

# Row 35
Set-TextCell "G35" "10"

# Row 36
Set-TextCell "G36" "10"

# Row 37
Set-TextCell "G37" "10"

# Row 38
Set-TextCell "G38" "10"

# Row 39
Set-TextCell "G39" "10"

# Row 40
Set-TextCell "E40" "3.87%"
Set-TextCell "G40" "10"

# Row 41
Set-TextCell "D41" "0.005643"
Set-TextCell "E41" "-6.65%"
Set-TextCell "G41" "10"

# Row 42
Set-TextCell "E42" "1.81%"
Set-TextCell "G42" "10"

# Row 43
Set-TextCell "D43" "0.002198"
Set-TextCell "E43" "-14.78%"
Set-TextCell "G43" "10"

# Row 44
Set-TextCell "D44" "0.01000"
Set-TextCell "E44" "22.35%"
Set-TextCell "G44" "10"

# Row 45
Set-TextCell "D45" "0.00005112"
Set-TextCell "E45" "-3.80%"
Set-TextCell "G45" "10"

# Row 46
Set-TextCell "E46" "-0.05%"
Set-TextCell "G46" "10"

# Row 47
Set-TextCell "D47" "0.08893"
Set-TextCell "G47" "10"

# Row 48
Set-TextCell "E48" "5.08%"
Set-TextCell "G48" "10"

# Row 49
Set-TextCell "D49" "0.00002098"
Set-TextCell "E49" "-0.05%"
Set-TextCell "G49" "10"

# Row 50
Set-TextCell "D50" "0.0001998"
Set-TextCell "E50" "-0.05%"
Set-TextCell "G50" "10"

# Row 51
Set-TextCell "G51" "10"
